$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 16, shifting existing rows 16:23 down to 17:24
$ws.Rows.Item(16).Insert()

# Populate the new weekly data row (row 16) with the same layout as its peers
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44582
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = 100112010
$ws.Cells.Item(16, 7).Value = "Achicoria"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 52
$ws.Cells.Item(16, 11).Value = 7000
$ws.Cells.Item(16, 12).Value = 7000
$ws.Cells.Item(16, 13).Value = 7000
$ws.Cells.Item(16, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(16, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 16).Value = 438
$ws.Cells.Item(16, 17).Value = 16
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Match the date-style formatting used by the sibling rows in column D
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
